$wb = $excel.ActiveWorkbook

# Map of cell -> new value for column F ("想去人数") that needs updating.
$updates = @{
    "F2"  = 334
    "F4"  = 10497
    "F6"  = 953
    "F7"  = 55
    "F9"  = 7356
    "F11" = 450
    "F13" = 132
    "F14" = 3221
    "F15" = 38
    "F17" = 694
    "F19" = 1042
    "F20" = 279
    "F21" = 81
    "F22" = 1668
}

# Both the "展览" sheet and the "全部类型" sheet carry the same rows and
# need the identical set of updates applied.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($cellRef in $updates.Keys) {
        $ws.Range($cellRef).Value = $updates[$cellRef]
    }
}
